$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 38084
$ws.Range("F4").Value = 646
$ws.Range("F5").Value = 807
$ws.Range("F6").Value = 494
$ws.Range("F7").Value = 384
$ws.Range("F8").Value = 472
$ws.Range("F9").Value = 878
$ws.Range("F11").Value = 789
$ws.Range("F12").Value = 614
$ws.Range("F13").Value = 99
$ws.Range("F15").Value = 50
$ws.Range("F16").Value = 705
$ws.Range("F17").Value = 199
$ws.Range("F18").Value = 504
$ws.Range("F20").Value = 1209
$ws.Range("F22").Value = 906
$ws.Range("F23").Value = 2639
$ws.Range("F24").Value = 1129
$ws.Range("F25").Value = 596
$ws.Range("F26").Value = 135
$ws.Range("F27").Value = 1192
$ws.Range("F28").Value = 47
$ws.Range("F29").Value = 872
$ws.Range("F30").Value = 84
$ws.Range("F31").Value = 1211

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 473

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 688

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 688
$ws.Range("F3").Value = 38084
$ws.Range("F5").Value = 646
$ws.Range("F6").Value = 807
$ws.Range("F7").Value = 494
$ws.Range("F9").Value = 384
$ws.Range("F10").Value = 472
$ws.Range("F11").Value = 473
$ws.Range("F15").Value = 878
$ws.Range("F17").Value = 789
$ws.Range("F18").Value = 614
$ws.Range("F19").Value = 99
$ws.Range("F25").Value = 50
$ws.Range("F27").Value = 705
$ws.Range("F28").Value = 199
$ws.Range("F29").Value = 504
$ws.Range("F31").Value = 1209
$ws.Range("F33").Value = 906
$ws.Range("F34").Value = 2639
$ws.Range("F35").Value = 1129
$ws.Range("F36").Value = 596
$ws.Range("F37").Value = 135
$ws.Range("F38").Value = 1192
$ws.Range("F39").Value = 47
$ws.Range("F41").Value = 872
$ws.Range("F42").Value = 84
$ws.Range("F43").Value = 1211
